{"js": "// Replace the 25 division-problem strings in the worksheet table, in\n// document order. Several original strings repeat (e.g. \"47\u00f75=\" occurs\n// twice, and \"74\u00f76=\"/\"30\u00f78=\" form a chain with the new values), so the\n// mapping is applied positionally \u2014 each occurrence, in document order,\n// is replaced by the corresponding entry in `replacements` \u2014 rather than\n// via a naive global text substitution.\nconst replacements = [\n  \"53\u00f76=\", \"82\u00f74=\", \"74\u00f77=\", \"79\u00f74=\", \"16\u00f74=\",\n  \"33\u00f72=\", \"31\u00f76=\", \"26\u00f77=\", \"25\u00f73=\", \"40\u00f74=\",\n  \"61\u00f72=\", \"79\u00f78=\", \"67\u00f75=\", \"72\u00f75=\", \"57\u00f78=\",\n  \"92\u00f74=\", \"48\u00f76=\", \"17\u00f73=\", \"61\u00f75=\", \"83\u00f73=\",\n  \"22\u00f74=\", \"28\u00f74=\", \"93\u00f78=\", \"30\u00f78=\", \"83\u00f74=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Collect every cell's first paragraph so we can read + rewrite its text\n// while walking cells in strict document (row-major) order.\nconst paragraphs = [];\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    cell.body.paragraphs.load(\"items,text\");\n    paragraphs.push(cell.body.paragraphs);\n  }\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const paraCollection of paragraphs) {\n  const para = paraCollection.items[0];\n  if (para && /\\d+\u00f7\\d+=/.test(para.text)) {\n    para.insertText(replacements[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem strings in the worksheet table, in\n# document order. Several original strings repeat (e.g. \"47\u00f75=\" occurs\n# twice, and \"74\u00f76=\"/\"30\u00f78=\" form a chain with the new values), so the\n# mapping is applied positionally -- each occurrence, in document order,\n# is replaced by the corresponding entry in $replacements -- rather than\n# via a naive global text substitution.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  \"53\u00f76=\", \"82\u00f74=\", \"74\u00f77=\", \"79\u00f74=\", \"16\u00f74=\",\n  \"33\u00f72=\", \"31\u00f76=\", \"26\u00f77=\", \"25\u00f73=\", \"40\u00f74=\",\n  \"61\u00f72=\", \"79\u00f78=\", \"67\u00f75=\", \"72\u00f75=\", \"57\u00f78=\",\n  \"92\u00f74=\", \"48\u00f76=\", \"17\u00f73=\", \"61\u00f75=\", \"83\u00f73=\",\n  \"22\u00f74=\", \"28\u00f74=\", \"93\u00f78=\", \"30\u00f78=\", \"83\u00f74=\"\n)\n\n$idx = 0\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $rng = $cell.Range\n    $rng.MoveEnd(1, -1) | Out-Null   # drop the cell-end mark, keep formatting\n    if ($rng.Text -match \"^[0-9]+.[0-9]+=$\") {\n      $rng.Text = $replacements[$idx]\n      $idx = $idx + 1\n    }\n  }\n}\n"}
